# Re-ran analysis.py on new synthetic volumes to correct angle calcs.
#
# The angle_min_deg/angle_max_deg/angle_mean_deg/angle_std_deg columns
# (T:W) on the "Synthetic" sheet were recomputed against the new
# synthetic-volume inputs; this writes the corrected values back in
# place. View/selection state on both sheets is also refreshed to match
# where the analyst last clicked after reviewing the results.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "SegThy"
$ws2 = $wb.Worksheets.Item(2)   # "Synthetic"

# --- Corrected angle statistics on "Synthetic" ---
$ws2.Range("T2").Value = -112.483254668527
$ws2.Range("U2").Value = 32.444948829019403
$ws2.Range("V2").Value = -92.680803784956595
$ws2.Range("W2").Value = 23.170010867658998
$ws2.Range("T3").Value = -74.9592337686245
$ws2.Range("U3").Value = 65.920813435909295
$ws2.Range("V3").Value = -56.083014187529599
$ws2.Range("W3").Value = 22.656017964062301
$ws2.Range("T4").Value = -22.515745623197699
$ws2.Range("U4").Value = 50.076068502654699
$ws2.Range("V4").Value = -5.66130044265776
$ws2.Range("W4").Value = 17.727140255431401
$ws2.Range("T5").Value = -106.688249581566
$ws2.Range("U5").Value = 41.573561062288903
$ws2.Range("V5").Value = -85.981800133711204
$ws2.Range("W5").Value = 26.553437563913501
$ws2.Range("T6").Value = -66.088725975100999
$ws2.Range("U6").Value = 10.640237642201299
$ws2.Range("V6").Value = -48.469623736432503
$ws2.Range("W6").Value = 18.771293550193299
$ws2.Range("T7").Value = -31.2275811003014
$ws2.Range("U7").Value = 37.498054018220699
$ws2.Range("V7").Value = -12.854720476895899
$ws2.Range("W7").Value = 20.0117022595701
$ws2.Range("T8").Value = -178.989102095037
$ws2.Range("U8").Value = 178.79866144223899
$ws2.Range("V8").Value = 18.1793537082957
$ws2.Range("W8").Value = 157.434413673243
$ws2.Range("T9").Value = -173.81513626117999
$ws2.Range("U9").Value = 170.07074452297499
$ws2.Range("V9").Value = -100.65642018092301
$ws2.Range("W9").Value = 44.094801302008101
$ws2.Range("T10").Value = -125.928984691991
$ws2.Range("U10").Value = -27.234624930914698
$ws2.Range("V10").Value = -63.0796939554129
$ws2.Range("T11").Value = -28.738070977922
$ws2.Range("U11").Value = 44.800276363492699
$ws2.Range("V11").Value = -10.8129706836328
$ws2.Range("W11").Value = 20.822539013200601
$ws2.Range("T12").Value = -120.855475984585
$ws2.Range("U12").Value = -54.222599557732799
$ws2.Range("V12").Value = -100.395944150444
$ws2.Range("T13").Value = -73.315389667923398
$ws2.Range("U13").Value = 73.617851897805593
$ws2.Range("V13").Value = -49.977564265835902
$ws2.Range("W13").Value = 28.377580218663699
$ws2.Range("T14").Value = -18.835588482217702
$ws2.Range("U14").Value = 53.382182517231499
$ws2.Range("V14").Value = -3.1006022436549601
$ws2.Range("W14").Value = 17.2440059382077
$ws2.Range("T15").Value = -117.76991223439001
$ws2.Range("U15").Value = -55.947542344323601
$ws2.Range("V15").Value = -96.244533959798304
$ws2.Range("T16").Value = -84.521341266425196
$ws2.Range("U16").Value = 94.827562968532305
$ws2.Range("V16").Value = -50.5049069866811
$ws2.Range("W16").Value = 37.8906191589803
$ws2.Range("T17").Value = -39.136416276149397
$ws2.Range("U17").Value = 52.599191885376499
$ws2.Range("V17").Value = -12.287104710591001
$ws2.Range("W17").Value = 27.399408084431901

# --- Refresh view/selection state ---
# "Synthetic": selection moves to J30; this also drops the stale
# horizontal scroll position (topLeftCell="M1") now that the selection
# sits back near the left edge of the sheet.
$ws2.Activate()
$ws2.Range("J30").Select()

# "SegThy" remains the active/tab-selected sheet, with its selection
# moved to V37.
$ws1.Activate()
$ws1.Range("V37").Select()
